$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2630.1484702492698, 2134.1724774265472, 2059.154798274425),
    @(2423.6079349026991, 1406.421030295031, 1701.2315806727984),
    @(2891.7262547040496, 2193.7536547415039, 2138.887987011994),
    @(2715.7902273637637, 2202.8567011645832, 2210.8064768191948),
    @(2940.5306950164641, 2148.034798690579, 2198.8020464359815),
    @(2785.7168650009248, 2303.563896204485, 2387.5631075826905),
    @(2397.5815275200744, 2111.8725650676633, 2071.2898662971816),
    @(2882.0547492517726, 2506.2593653705144, 2322.0381252657562),
    @(3159.4912255857453, 2252.6815124387631, 2035.8815171880751),
    @(2442.5676874809978, 1443.2985593610731, 1578.0775832209015),
    @(2285.2872783357448, 1670.3501514908894, 1611.6230629012018),
    @(3074.4885482261579, 2504.9068044909186, 2362.1270015828272),
    @(3145.5104891155056, 2517.3901677265812, 2258.3211231646919),
    @(3159.7893428970565, 2565.8316250216185, 2289.3730090404611),
    @(3128.9397263965379, 2515.7569156154477, 2389.3836037915657),
    @(3012.36363244525, 2161.4227746618949, 1838.1326222642756)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
